$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header labels (Volume/Number and week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 56
$ws.Range("J15").Value = 44
$ws.Range("K15").Value = 27.272727272727
$ws.Range("L15").Value = 107.407407407407
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = -11.111111111111
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 31.25
$ws.Range("I16").Value = 213
$ws.Range("J16").Value = 212
$ws.Range("K16").Value = 0.471698113207
$ws.Range("L16").Value = -8.974358974358
$ws.Range("M16").Value = -27.551020408163
$ws.Range("N16").Value = -77.219251336898
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -5.882352941176
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = -7.246376811594
$ws.Range("I17").Value = 728
$ws.Range("J17").Value = 682
$ws.Range("K17").Value = 6.74486803519
$ws.Range("L17").Value = 6.432748538011
$ws.Range("M17").Value = 112.244897959184
$ws.Range("N17").Value = -15.74074074074
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 217
$ws.Range("J18").Value = 234
$ws.Range("K18").Value = -7.264957264957
$ws.Range("L18").Value = -12.85140562249
$ws.Range("M18").Value = -49.061032863849
$ws.Range("N18").Value = -91.106557377049
$ws.Range("C19").Value = 40
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 143
$ws.Range("G19").Value = 107
$ws.Range("H19").Value = 33.644859813084
$ws.Range("I19").Value = 1040
$ws.Range("J19").Value = 1096
$ws.Range("K19").Value = -5.109489051094
$ws.Range("L19").Value = -8.208296557811
$ws.Range("M19").Value = 52.716593245227
$ws.Range("N19").Value = -8.851884312007
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -55.555555555555
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = -36.666666666666
$ws.Range("I20").Value = 162
$ws.Range("J20").Value = 252
$ws.Range("K20").Value = -35.714285714285
$ws.Range("L20").Value = -50.609756097561
$ws.Range("M20").Value = -34.677419354838
$ws.Range("N20").Value = -95.24787327662
$ws.Range("D21").Value = 68
$ws.Range("E21").Value = 0
$ws.Range("G21").Value = 261
$ws.Range("H21").Value = 4.597701149425
$ws.Range("I21").Value = 2418
$ws.Range("J21").Value = 2525
$ws.Range("K21").Value = -4.237623762376
$ws.Range("L21").Value = -9.472107824784
$ws.Range("M21").Value = 18.471337579617
$ws.Range("N21").Value = -72.739571589628
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = -40
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = -70.588235294117
$ws.Range("I23").Value = 106
$ws.Range("J23").Value = 93
$ws.Range("K23").Value = 13.978494623655
$ws.Range("L23").Value = 6
$ws.Range("M23").Value = 112
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 69
$ws.Range("E24").Value = -5.797101449275
$ws.Range("F24").Value = 286
$ws.Range("G24").Value = 280
$ws.Range("H24").Value = 2.142857142857
$ws.Range("I24").Value = 2877
$ws.Range("J24").Value = 2857
$ws.Range("K24").Value = 0.70003500175
$ws.Range("L24").Value = -4.418604651162
$ws.Range("M24").Value = 6.674082313681
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 41
$ws.Range("E25").Value = -9.756097560975
$ws.Range("F25").Value = 157
$ws.Range("G25").Value = 151
$ws.Range("H25").Value = 3.973509933774
$ws.Range("I25").Value = 1677
$ws.Range("J25").Value = 1517
$ws.Range("K25").Value = 10.547132498352
$ws.Range("L25").Value = 18.851878100637
$ws.Range("C26").Value = 41
$ws.Range("E26").Value = 5.128205128205
$ws.Range("F26").Value = 160
$ws.Range("G26").Value = 158
$ws.Range("H26").Value = 1.26582278481
$ws.Range("I26").Value = 1355
$ws.Range("J26").Value = 1347
$ws.Range("K26").Value = 0.593912397921
$ws.Range("L26").Value = 7.284243863816
$ws.Range("M26").Value = -5.17844646606
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 73
$ws.Range("K27").Value = -2.739726027397
$ws.Range("L27").Value = 54.347826086956
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = -21.052631578947
$ws.Range("I28").Value = 155
$ws.Range("J28").Value = 147
$ws.Range("K28").Value = 5.442176870748
$ws.Range("L28").Value = 0.64935064935
$ws.Range("G31").Value = 3
$ws.Range("J31").Value = 17
$ws.Range("K31").Value = -35.294117647058
$ws.Range("L31").Value = 0
$ws.Range("L33").Value = -61.538461538461

# --- Cells changing from numeric to "N/A" text style (copy format from an existing text cell, e.g. C14) ---
$ws.Range("G14").Value = "0"
$ws.Range("H14").Value = "***.*"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("G14","H14","C29","D29","E29","C30","D30","E30","D33","E33").PasteSpecial(-4122)

# --- Cells changing from "N/A" text to numeric (copy format from existing numeric cells with matching format) ---
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("C15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E15").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Applied all crime data updates"
